# "Mise a jour de la liste des taches" :
#  - plusieurs taches marquees "A faire" (colonne B) passent soit a
#    "En cours" (colonne C) soit a "Termine" (colonne D), avec mise a jour
#    du responsable (colonne E) et des dates de debut/fin (colonnes F/G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

function Copy-DateFormat($Row, $Col) {
    # Clone le format (deja numerique/date) d'une cellule de reference qui
    # porte deja la mise en forme "date" voulue, plutot que de creer un
    # nouveau style lors de l'ecriture de la date.
    $ws.Range("F37").Copy() | Out-Null
    $ws.Cells.Item($Row, $Col).PasteSpecial(-4122) | Out-Null
}

function Set-TacheTerminee($Row, $Responsable, $DateDebut, $DateFin) {
    # La tache n'est plus "A faire" -> on vide la croix en colonne B.
    $ws.Cells.Item($Row, 2).Value = ""

    # Elle est desormais "Terminee" -> la croix passe en colonne D.
    $ws.Cells.Item($Row, 4).Value = "X"

    # Responsable de la tache.
    $ws.Cells.Item($Row, 5).Value = $Responsable

    # Dates de debut (F) et de fin (G).
    Copy-DateFormat $Row 6
    $ws.Cells.Item($Row, 6).Value2 = $DateDebut

    Copy-DateFormat $Row 7
    $ws.Cells.Item($Row, 7).Value2 = $DateFin

    $excel.CutCopyMode = 0
}

function Set-TacheEnCours($Row, $Responsable) {
    # La tache n'est plus "A faire" -> on vide la croix en colonne B.
    $ws.Cells.Item($Row, 2).Value = ""

    # Elle est desormais "En cours" -> la croix passe en colonne C.
    $ws.Cells.Item($Row, 3).Value = "X"

    # Responsable de la tache.
    $ws.Cells.Item($Row, 5).Value = $Responsable
}

function Set-TacheEnCoursAvecDebut($Row, $Responsable, $DateDebut) {
    Set-TacheEnCours $Row $Responsable

    Copy-DateFormat $Row 6
    $ws.Cells.Item($Row, 6).Value2 = $DateDebut

    $excel.CutCopyMode = 0
}

# Serial Excel des dates (base 30/12/1899) : 28/09/2016 -> 42641, 27/09/2016 -> 42640.
$dateDebutFin28 = 42641
$date27 = 42640

# Taches passees "Terminee", responsable Killian, du 28/09/2016 au 28/09/2016.
Set-TacheTerminee 36 "Killian" $dateDebutFin28 $dateDebutFin28
Set-TacheTerminee 37 "Killian" $dateDebutFin28 $dateDebutFin28
Set-TacheTerminee 68 "Killian" $dateDebutFin28 $dateDebutFin28
Set-TacheTerminee 69 "Killian" $dateDebutFin28 $dateDebutFin28

# Taches passees "En cours", responsable Tony (nouveau), sans dates.
Set-TacheEnCours 50 "Tony"
Set-TacheEnCours 51 "Tony"

# Taches passees "En cours", responsable Erwann, avec une date de debut.
Set-TacheEnCoursAvecDebut 62 "Erwann" $date27
Set-TacheEnCoursAvecDebut 63 "Erwann" $date27

# Repositionne la vue sur la feuille comme dans le classeur d'origine.
$ws.Range("I64").Select() | Out-Null
